$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new columns (G:K) for the new meta metrics, shifting the
# existing arrecadado_* ... maior_ano columns from G:V to L:AA.
$ws.Range("G1:K1").EntireColumn.Insert()

# Header row (row 1): meta, meta_avg, meta_std, meta_min, meta_max
$ws.Cells.Item(1, 7).Value = "meta"
$ws.Cells.Item(1, 8).Value = "meta_avg"
$ws.Cells.Item(1, 9).Value = "meta_std"
$ws.Cells.Item(1, 10).Value = "meta_min"
$ws.Cells.Item(1, 11).Value = "meta_max"

# Give the 5 new data columns the same currency format as the other
# style-3 columns (arrecadado_*, apoio_*); only the data rows (2-6),
# so the header row keeps its original bold/border style.
$ws.Range("G2:K6").NumberFormat = "R$ #,##0.00"

# Data values for the 5 new meta columns, row by row.
$ws.Cells.Item(2, 7).Value = 721610.3061912227
$ws.Cells.Item(2, 8).Value = 10458.12037958294
$ws.Cells.Item(2, 9).Value = 11144.2267578863
$ws.Cells.Item(2, 10).Value = 44.33046360042423
$ws.Cells.Item(2, 11).Value = 50590.198657868

$ws.Cells.Item(3, 7).Value = 5883940.636230236
$ws.Cells.Item(3, 8).Value = 13372.59235506872
$ws.Cells.Item(3, 9).Value = 19267.96260047285
$ws.Cells.Item(3, 10).Value = 23.98859826184044
$ws.Cells.Item(3, 11).Value = 147790.8327903106

$ws.Cells.Item(4, 7).Value = 1712986.472842461
$ws.Cells.Item(4, 8).Value = 9732.877686604894
$ws.Cells.Item(4, 9).Value = 10102.88946115519
$ws.Cells.Item(4, 10).Value = 46.55761904502517
$ws.Cells.Item(4, 11).Value = 83151.82469725677

$ws.Cells.Item(5, 7).Value = 7150010.825257363
$ws.Cells.Item(5, 8).Value = 10347.33838676898
$ws.Cells.Item(5, 9).Value = 16064.05218382809
$ws.Cells.Item(5, 10).Value = 12.04441558726698
$ws.Cells.Item(5, 11).Value = 198811.9434626772

$ws.Cells.Item(6, 7).Value = 131168.4623975197
$ws.Cells.Item(6, 8).Value = 18738.35177107424
$ws.Cells.Item(6, 9).Value = 19781.31029827062
$ws.Cells.Item(6, 10).Value = 2420.445520432476
$ws.Cells.Item(6, 11).Value = 54319.48382898097

Write-Host "Applied meta columns."
